$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, shifting existing rows 38:81 down to 39:82
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with this week's data
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 44651
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = 100112031
$ws.Cells.Item(38, 7).Value = "Poroto verde"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 80
$ws.Cells.Item(38, 11).Value = 24000
$ws.Cells.Item(38, 12).Value = 25000
$ws.Cells.Item(38, 13).Value = 24500
$ws.Cells.Item(38, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Región del Maule"
$ws.Cells.Item(38, 16).Value = 980
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# Match the date cell style used by the other date cells in column D
$ws.Range("D39").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$excel.CutCopyMode = 0
